$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers for columns B and C
$ws.Range("B1").Value = "CategoryName"
$ws.Range("C1").Value = "UserName"

# Copy the formatting already used by column C (Consolas font, centered)
# onto column B so the new CategoryName column matches UserName's look.
$ws.Range("C2:C11").Copy()
$ws.Range("B2:B11").PasteSpecial(-4122)  # xlPasteFormats

# Fill the data rows with the new constant values
$ws.Range("B2:B11").Value = "Thực tập"
$ws.Range("C2:C11").Value = "Nguyendinhnam28803"

# Restore the active selection
$ws.Range("H14").Select()
